$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue([string]$cellRef, [string]$val) {
    # Forces a numeric-looking string to be stored as text, matching the
    # original cell type, without leaving a NumberFormat/quotePrefix style
    # behind: write with a leading apostrophe (forces text), then copy/paste
    # the format (only) from a known plain-text cell (D2) on top of it so
    # no extra style index is introduced.
    $cell = $ws.Range($cellRef)
    $cell.Value = "'" + $val
    $ws.Range("D2").Copy()
    $cell.PasteSpecial(-4122)
}

$ws.Range("D2").Value = "26.080.86"
$ws.Range("E2").Value = "  +3.26%  "
$ws.Range("D3").Value = "1.603.62"
$ws.Range("E3").Value = "  +3.37%  "
$ws.Range("E4").Value = "  -0.14%  "
Set-TextValue "D5" "212.85"
$ws.Range("E5").Value = "  +2.89%  "
$ws.Range("E6").Value = "  -0.17%  "
Set-TextValue "D7" "0.487"
$ws.Range("E7").Value = "  +1.99%  "
$ws.Range("E8").Value = "  +2.32%  "
$ws.Range("E9").Value = "  +0.94%  "
Set-TextValue "D10" "18.05"
$ws.Range("E10").Value = "  +2.15%  "
Set-TextValue "D11" "0.0817"
$ws.Range("E11").Value = "  +4.66%  "
$ws.Range("D12").Value = "1.827.05"
$ws.Range("E12").Value = "  +3.40%  "
$ws.Range("D13").Value = "1.607.90"
$ws.Range("E13").Value = "  +4.50%  "
$ws.Range("E14").Value = "  +0.52%  "
Set-TextValue "D15" "0.511"
$ws.Range("E15").Value = "  +1.44%  "
$ws.Range("D16").Value = "26.092.88"
$ws.Range("E16").Value = "  +3.35%  "
Set-TextValue "D17" "60.38"
$ws.Range("E17").Value = "  +3.08%  "
$ws.Range("D18").Value = "0.0₃0722"
$ws.Range("E18").Value = "  +2.01%  "
$ws.Range("E19").Value = "  -0.18%  "
Set-TextValue "D20" "201.72"
$ws.Range("E20").Value = "  +8.83%  "
Set-TextValue "D21" "4.22"
$ws.Range("E21").Value = "  +2.97%  "
$ws.Range("E22").Value = "  +0.39%  "
Set-TextValue "D23" "5.99"
$ws.Range("E23").Value = "  +2.57%  "
Set-TextValue "D24" "1.86"
$ws.Range("E24").Value = "  +13.97%  "
Set-TextValue "D25" "141.63"
$ws.Range("E25").Value = "  +1.62%  "
$ws.Range("E26").Value = "  -0.15%  "
$ws.Range("E27").Value = "  -4.39%  "
$ws.Range("E28").Value = "  +2.31%  "
$ws.Range("E29").Value = "  +0.73%  "
$ws.Range("E30").Value = "  +2.01%  "
Set-TextValue "D31" "0.0473"
$ws.Range("E31").Value = "  +1.60%  "
$ws.Range("E32").Value = "  +2.82%  "
$ws.Range("E33").Value = "  -0.21%  "
$ws.Range("E34").Value = "  +1.64%  "
$ws.Range("E35").Value = "  +0.90%  "
$ws.Range("D36").Value = "1.124.84"
$ws.Range("E36").Value = "  +3.93%  "
$ws.Range("E37").Value = "  +9.84%  "
$ws.Range("E38").Value = "  -0.07%  "
$ws.Range("B39").Value = "MXToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue "D39" "2.30"
$ws.Range("E39").Value = "  +2.51%  "
$ws.Range("B40").Value = "ARBITRUM"
$ws.Range("C40").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextValue "D40" "0.782"
$ws.Range("E40").Value = "  +2.57%  "
Set-TextValue "D41" "0.492"
$ws.Range("E41").Value = "  -0.37%  "
$ws.Range("E43").Value = "  +2.52%  "
$ws.Range("D44").Value = "1.739.62"
$ws.Range("E44").Value = "  +3.44%  "
Set-TextValue "D45" "92.75"
$ws.Range("E45").Value = "  +0.25%  "
$ws.Range("E46").Value = "  +4.22%  "
Set-TextValue "D47" "53.49"
$ws.Range("E47").Value = "  +2.36%  "
Set-TextValue "D48" "0.0505"
$ws.Range("E48").Value = "  +0.50%  "
Set-TextValue "D49" "0.408"
$ws.Range("E49").Value = "  +1.11%  "
$ws.Range("E50").Value = "  +0.05%  "
$ws.Range("D51").Value = "0.0₇0946"
$ws.Range("E51").Value = "  -14.65%  "

$excel.CutCopyMode = 0
